# Update the "取得日時" (acquired datetime) column to reflect the latest
# scrape run timestamp: 2025-12-28 18:34:21 (JST), replacing the previous
# 2025-12-28 18:26:33 value for every data row on the "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-28 18:34:21"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
